# Fix typo "tu[eDecl" -> "typeDecl" in the third code example on the
# "Simple Grammar Examples" slide (slide 21).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange

# Paragraph 6 reads:  "  initialDecl = tu[eDecl ."
$para = $tr.Paragraphs(6, 1)

# Locate "tu[eDecl" (8 characters) within the paragraph and replace it
# with "typeDecl", leaving the surrounding text/runs untouched.
$typoStart = $para.Text.IndexOf("tu[eDecl") + 1
$target = $para.Characters($typoStart, 8)
$target.Text = "typeDecl"
